$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text string into a cell without Excel's automatic
# type coercion (e.g. turning "Feb 2026" into a date) and without touching
# the cell's style. We do this by placing a string-literal formula in the
# cell, then copying the cell and pasting back "values only" - the pasted
# value retains its string type, and no NumberFormat/style gets applied.
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $formulaText = '="' + $escaped + '"'
    $cell.Formula = $formulaText
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = $false

# Row data: row number -> month/year text (short form) and timestamp text.
$rows = @(2, 3, 4, 5, 6, 7)
$months = @("Feb 2026", "Mar 2026", "Sep 2026", "Nov 2026", "Nov 2026", "Jan 2027")
$stamps = @("2026-02-16 11:06:05", "2026-02-16 11:06:07", "2026-02-16 11:06:07", "2026-02-16 11:06:07", "2026-02-16 11:06:07", "2026-02-16 11:06:08")

# First pass: populate column C with the full "Expected Launch : <month>"
# text for every row. This is what originally seeded those longer strings
# into the workbook's shared string table.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $m = $months[$i]
    $longText = "Expected Launch : " + $m
    $cCell = $ws.Cells.Item($r, 3)
    Set-TextValue $cCell $longText
}

# Second pass: overwrite column C with just the short month/year text, and
# fill column D with the timestamp text.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $m = $months[$i]
    $s = $stamps[$i]

    $cCell = $ws.Cells.Item($r, 3)
    Set-TextValue $cCell $m

    $dCell = $ws.Cells.Item($r, 4)
    Set-TextValue $dCell $s
}

$excel.CutCopyMode = $false
